$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

Set-TextValue "D2" "48.517.12"
Set-TextValue "E2" "  +7.50%  "

Set-TextValue "D3" "2.650.81"
Set-TextValue "E3" "  +11.00%  "

Set-TextValue "D4" "0.997"
Set-TextValue "E4" "  -0.39%  "

Set-TextValue "D5" "314.35"
Set-TextValue "E5" "  +7.46%  "

Set-TextValue "D6" "106.30"
Set-TextValue "E6" "  +13.82%  "

Set-TextValue "D7" "0.615"
Set-TextValue "E7" "  +10.66%  "

Set-TextValue "D8" "0.996"
Set-TextValue "E8" "  -0.42%  "

Set-TextValue "E9" "  +20.07%  "

Set-TextValue "D10" "40.84"
Set-TextValue "E10" "  +19.64%  "

Set-TextValue "D11" "0.0867"
Set-TextValue "E11" "  +11.95%  "

Set-TextValue "D12" "55.74"
Set-TextValue "E12" "  +4.47%  "

Set-TextValue "D13" "8.44"
Set-TextValue "E13" "  +21.03%  "

Set-TextValue "D14" "3.039.67"
Set-TextValue "E14" "  +10.41%  "

Set-TextValue "E15" "  +3.62%  "

Set-TextValue "D16" "2.651.52"
Set-TextValue "E16" "  +11.01%  "

Set-TextValue "D17" "0.947"
Set-TextValue "E17" "  +14.84%  "

Set-TextValue "D18" "15.46"
Set-TextValue "E18" "  +9.86%  "

Set-TextValue "D19" "48.358.44"
Set-TextValue "E19" "  +7.14%  "

Set-TextValue "E20" "  +11.59%  "

Set-TextValue "D21" "13.44"
Set-TextValue "E21" "  +8.64%  "

Set-TextValue "D22" "6.91"
Set-TextValue "E22" "  +13.67%  "

Set-TextValue "D23" "73.91"
Set-TextValue "E23" "  +11.30%  "

Set-TextValue "D24" "282.43"
Set-TextValue "E24" "  +18.99%  "

Set-TextValue "D25" "3.13"
Set-TextValue "E25" "  +13.56%  "

Set-TextValue "D26" "2.26"
Set-TextValue "E26" "  +19.73%  "

Set-TextValue "D27" "30.34"
Set-TextValue "E27" "  +45.50%  "

Set-TextValue "E28" "  +0.20%  "

Set-TextValue "D29" "4.13"
Set-TextValue "E29" "  +2.72%  "

Set-TextValue "D30" "10.80"
Set-TextValue "E30" "  +13.56%  "

Set-TextValue "B31" "InjectiveProtocol"
Set-TextValue "C31" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D31" "40.49"
Set-TextValue "E31" "  +8.60%  "

Set-TextValue "B32" "Toncoin"
Set-TextValue "C32" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D32" "2.32"
Set-TextValue "E32" "  +4.71%  "

Set-TextValue "D33" "6.27"
Set-TextValue "E33" "  +16.55%  "

Set-TextValue "E34" "  -2.20%  "

Set-TextValue "D35" "0.0867"
Set-TextValue "E35" "  +14.69%  "

Set-TextValue "D36" "2.90"
Set-TextValue "E36" "  +7.06%  "

Set-TextValue "D37" "2.26"
Set-TextValue "E37" "  +15.01%  "

Set-TextValue "D38" "153.41"
Set-TextValue "E38" "  +4.65%  "

Set-TextValue "D39" "0.127"
Set-TextValue "E39" "  +13.50%  "

Set-TextValue "E40" "  +10.07%  "

Set-TextValue "D41" "16.68"
Set-TextValue "E41" "  +15.59%  "

Set-TextValue "E42" "  +18.57%  "

Set-TextValue "D43" "23.09"
Set-TextValue "E43" "  +49.43%  "

Set-TextValue "D44" "3.77"
Set-TextValue "E44" "  +19.32%  "

Set-TextValue "D45" "0.0339"
Set-TextValue "E45" "  +16.45%  "

Set-TextValue "D46" "2.203.46"
Set-TextValue "E46" "  +11.81%  "

Set-TextValue "D47" "99.41"
Set-TextValue "E47" "  +12.83%  "

Set-TextValue "D48" "0.998"
Set-TextValue "E48" "  -0.07%  "

Set-TextValue "D49" "9.97"
Set-TextValue "E49" "  +18.81%  "

Set-TextValue "B50" "Stacks"
Set-TextValue "C50" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D50" "1.91"
Set-TextValue "E50" "  +12.98%  "

Set-TextValue "B51" "Aave"
Set-TextValue "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D51" "115.47"
Set-TextValue "E51" "  +16.59%  "
